$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for columns B and C (rows 2-13); column D values are unchanged.
$newB = @{2=67.5; 3=83; 4=70.75; 5=72.75; 6=73.75; 7=71.25; 8=68.75; 9=74.75; 10=70.75; 11=74.5; 12=73; 13=72.25}
$newC = @{2=62.5; 3=82.5; 4=62.5; 5=52.5; 6=45; 7=47.5; 8=42.5; 9=40; 10=38.75; 11=36; 12=32.5; 13=30}

for ($r = 2; $r -le 13; $r++) {
    $ws.Range("B$r").Value = $newB[$r]
    $ws.Range("C$r").Value = $newC[$r]
    $ws.Range("E$r").Formula = "=B$r-C$r"
}

# Columns B, C and E (rows 2-13) get a centered integer number format.
$ws.Range("B2:C13").NumberFormat = "0"
$ws.Range("B2:C13").HorizontalAlignment = -4108
$ws.Range("E2:E13").NumberFormat = "0"
$ws.Range("E2:E13").HorizontalAlignment = -4108

# Column D (rows 2-13) keeps its General number format but becomes centered.
$ws.Range("D2:D13").HorizontalAlignment = -4108

# Move the active selection to B4.
$ws.Range("B4").Select()
